$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Switch to the "Repayment schedule" tab (was "Summary") and select it.
$ws.Activate()

# Insert a new blank column before column N, shifting the old N/O/P
# (Late / heading / Outstanding) data right to O/P/Q.
$ws.Columns("N").Insert()

# The freshly inserted column N picks up the width of its left neighbour
# (column M), matching what Excel does on a manual column insert.
$ws.Columns("N").ColumnWidth = 9.83

# Leave the selection on the newly active sheet at R7.
$ws.Range("R7").Select() | Out-Null
